# Update workbook to the new scrape snapshot.
# Two sheets are affected: "展览" (insert point row 18) and "全部类型" (insert point row 19).
# Both get: one new row inserted with a new event, and a handful of F-column
# (想去人数) value bumps, some on rows before the insertion point and some on
# rows that end up shifted down by the insertion.

$wb = $excel.ActiveWorkbook

function Add-Event($ws, $insertRow, $date, $name, $place, $timeRange, $wantCount, $price, $link, $cover) {
    # Insert a brand new row, shifting everything at/after $insertRow down by one.
    $ws.Rows.Item($insertRow).Insert()

    # Give the new row's index cell (column A) the same look (bold / bordered /
    # centered) as the rest of the index column by copying format from the row
    # that is now directly above the freshly inserted row.
    $ws.Range("A" + ($insertRow - 1)).Copy()
    $ws.Range("A" + $insertRow).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($insertRow, 1).Value = $insertRow - 1

    # The date column holds plain text like "2024-10-01" in every other row
    # (not a real Date value). Force text entry with a leading apostrophe,
    # then strip the resulting quote-prefix style so the cell matches its
    # siblings (no style attribute at all).
    $ws.Cells.Item($insertRow, 2).Value = "'" + $date
    $ws.Cells.Item($insertRow, 2).ClearFormats()

    $ws.Cells.Item($insertRow, 3).Value = $name
    $ws.Cells.Item($insertRow, 4).Value = $place
    $ws.Cells.Item($insertRow, 5).Value = $timeRange
    $ws.Cells.Item($insertRow, 6).Value = $wantCount
    $ws.Cells.Item($insertRow, 7).Value = $price
    $ws.Cells.Item($insertRow, 8).Value = $link
    $ws.Cells.Item($insertRow, 9).Value = $cover
}

function Set-Count($ws, $row, $value) {
    $ws.Cells.Item($row, 6).Value = $value
}

function Append-IndexRow($ws) {
    $lastRow = $ws.UsedRange.Rows.Count
    $ws.Range("A" + ($lastRow - 1)).Copy()
    $ws.Range("A" + $lastRow).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $ws.Cells.Item($lastRow, 1).Value = $lastRow - 1
}

# Column A is a plain 0-based row index (A(n) = n - 1) for every single data
# row. When Excel physically inserts a row, the cells below the insertion
# point keep travelling with their old A value, so column A has to be
# re-stamped for every row from the insertion point down to the (new) last
# row to keep it a simple position-based sequence.
function Renumber-IndexColumn($ws, $fromRow) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = $fromRow; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

function Update-Sheet($ws, $insertRow, $preBumps, $postBumps) {
    # ---- F-column bumps on rows that sit before the insertion point ----
    foreach ($r in $preBumps.Keys) {
        Set-Count $ws $r $preBumps[$r]
    }

    # ---- Insert the brand-new event row ----
    Add-Event $ws $insertRow "2024-10-01" "赣州·十万伏特-第八届青年文化综合展览会" "东阳山路9-3号 温州街" "2024.10.01 10:00-10.02 17:00" 2 55 "https://show.bilibili.com/platform/detail.html?id=92405" "//i1.hdslb.com/bfs/openplatform/202409/OUu2pvUR1726127605875.jpeg"

    # ---- F-column bumps on rows that were shifted down by the insertion ----
    # Keys are the NEW row numbers (after the insert above).
    foreach ($r in $postBumps.Keys) {
        Set-Count $ws $r $postBumps[$r]
    }

    # ---- Append the brand new trailing index row (column A only; its B..I
    # content already exists because it was shifted down from the old last
    # row) ----
    Append-IndexRow $ws
}

# Pre-insertion F-column bumps, keyed by row -> new value.
$preBumps1 = @{ 2 = 7129; 3 = 21; 5 = 20; 7 = 171; 12 = 207; 13 = 6; 14 = 454; 16 = 1843 }
$preBumps4 = @{ 2 = 7129; 3 = 21; 5 = 20; 8 = 171; 13 = 207; 14 = 6; 15 = 454; 17 = 1843 }

# Post-insertion F-column bumps (row numbers are the NEW, post-insert row numbers).
$postBumps1 = @{ 19 = 3698; 23 = 32; 25 = 2344; 27 = 280; 28 = 12; 29 = 1; 30 = 38; 33 = 1; 34 = 21; 35 = 160; 36 = 1380; 37 = 125 }
$postBumps4 = @{ 20 = 3698; 24 = 32; 26 = 2344; 28 = 280; 29 = 12; 30 = 1; 31 = 38; 34 = 1; 35 = 21; 36 = 160; 37 = 1380; 38 = 125 }

$ws1 = $wb.Worksheets.Item("展览")
Update-Sheet $ws1 18 $preBumps1 $postBumps1

$ws4 = $wb.Worksheets.Item("全部类型")
Update-Sheet $ws4 19 $preBumps4 $postBumps4
